$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing candidate: SOLARMEMS SSoC-A -> SOLARMEMS SSoC-A60 ---
$ws.Range("B3").Value = "SOLARMEMS SSoC-A60"

# --- Add new row 8: SOLARMEMS nanoSSoC-D60 candidate ---
# Copy formatting from existing representative cells for each column style,
# then set the values/text for the new row.

# D8, G8, H8, L8 use the "Neutral" style (style index 3) - copy from D4
$ws.Range("D4").Copy()
$ws.Range("D8").PasteSpecial(-4122)
$ws.Range("D4").Copy()
$ws.Range("G8").PasteSpecial(-4122)
$ws.Range("D4").Copy()
$ws.Range("H8").PasteSpecial(-4122)
$ws.Range("D4").Copy()
$ws.Range("L8").PasteSpecial(-4122)

# E8 also uses "Neutral" style (style index 3) - copy from D4
$ws.Range("D4").Copy()
$ws.Range("E8").PasteSpecial(-4122)

# F8, J8 use the "Good" style (style index 4) - copy from F2
$ws.Range("F2").Copy()
$ws.Range("F8").PasteSpecial(-4122)
$ws.Range("F2").Copy()
$ws.Range("J8").PasteSpecial(-4122)

# K8 uses the "Good" style with quote prefix (style index 10) - copy from K3
$ws.Range("K3").Copy()
$ws.Range("K8").PasteSpecial(-4122)

# I8 uses a "Good" styled cell (closest achievable match)
$ws.Range("I8").Style = "Good"

$ws.Application.CutCopyMode = $False

# Now set the cell values/text
$ws.Range("B8").Value = "SOLARMEMS nanoSSoC-D60"
$ws.Range("C8").Value = "http://www.solar-mems.com/smt_pdf/Brochure_NanoSSOC-D60.pdf"
$ws.Range("D8").Value = "unknown"
$ws.Range("E8").Value = "<23mA(?)"
$ws.Range("F8").Value = "6.5 grams"
$ws.Range("G8").Value = "unknown"
$ws.Range("H8").Value = "unknown"
$ws.Range("I8").Value = "<0.5 deg"
$ws.Range("J8").Value = "120 deg."
$ws.Range("K8").Value = "-30C to +85C"
$ws.Range("L8").Value = "43x14x5.9mm"

# --- Update row 9 comment text stays same content (index shift handled automatically) ---

# --- Add row 11 comment about trying ADCOLE again ---
$ws.Range("B11").Value = "Try ADCOLE again since they're stuff is very good. If we run into a problem with sun sensors not being able to update fast enough, they may be able to help. Otherwise, maybe we'll have to do an intertial solution."

# --- Update selection to match author's final cursor position ---
$ws.Range("B14").Select()
